$d = $word.ActiveDocument

# The document starts with a single 2-column table holding "TEAM ID" and
# "PROJECT NAME" rows. We need to insert a new first row containing a
# "DATE" / "09-10-22" pair above the existing "TEAM ID" row.
$table = $d.Tables.Item(1)
$firstRow = $table.Rows.Item(1)

# Adding a row "before" the current first row inserts a brand-new row at
# the top of the table, inheriting the table's existing cell/paragraph
# formatting (so the new cells keep the sz=28/szCs=28 run formatting used
# throughout the rest of the table).
$newRow = $table.Rows.Add($firstRow)

$newRow.Cells.Item(1).Range.Text = "DATE"
$newRow.Cells.Item(2).Range.Text = "09-10-22"
